# Update Bankrate mortgage rate data (interest_rate / apr / updated_date)
# Values like "6.67%" and dates like "2025-08-11" must stay as literal text
# (matching the original shared-string layout), not get auto-converted by
# Excel into percentage/date numbers. Force text via NumberFormat "@" then
# reset the cell style back to Normal so no numeric style sticks around.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# row 2: 30-Year Fixed Rate
Set-TextValue $ws.Range("B2") "6.67%"
Set-TextValue $ws.Range("C2") "6.73%"
Set-TextValue $ws.Range("F2") "2025-08-11"

# row 3: 20-Year Fixed Rate
Set-TextValue $ws.Range("B3") "6.36%"
Set-TextValue $ws.Range("C3") "6.45%"
Set-TextValue $ws.Range("F3") "2025-08-11"

# row 4: 15-Year Fixed Rate
Set-TextValue $ws.Range("B4") "5.85%"
Set-TextValue $ws.Range("C4") "5.94%"
Set-TextValue $ws.Range("F4") "2025-08-11"

# row 5: 10-Year Fixed Rate
Set-TextValue $ws.Range("B5") "5.76%"
Set-TextValue $ws.Range("C5") "5.83%"
Set-TextValue $ws.Range("F5") "2025-08-11"

# row 6: 30-Year Fixed Rate FHA
Set-TextValue $ws.Range("B6") "6.56%"
Set-TextValue $ws.Range("C6") "6.61%"
Set-TextValue $ws.Range("F6") "2025-08-11"

# row 7: 30-Year Fixed Rate VA
Set-TextValue $ws.Range("B7") "6.66%"
Set-TextValue $ws.Range("C7") "6.70%"
Set-TextValue $ws.Range("F7") "2025-08-11"

# row 8: 30-Year Fixed Rate Jumbo
Set-TextValue $ws.Range("B8") "6.67%"
Set-TextValue $ws.Range("F8") "2025-08-11"
